$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 12 (old last data row); remaining rows keep their row numbers (2-11)
$ws.Rows.Item(12).Delete()

# Update remaining data rows (2-11) with the new values
$ws.Cells.Item(2, 2).Value = 0.3963
$ws.Cells.Item(2, 3).Value = 283
$ws.Cells.Item(2, 4).Value = 1575
$ws.Cells.Item(2, 5).Value = 0.00113
$ws.Cells.Item(3, 2).Value = 0.3711
$ws.Cells.Item(3, 3).Value = 71
$ws.Cells.Item(3, 4).Value = 358
$ws.Cells.Item(3, 5).Value = 0.000319
$ws.Cells.Item(4, 2).Value = 0.3512
$ws.Cells.Item(4, 3).Value = 45
$ws.Cells.Item(4, 4).Value = 223
$ws.Cells.Item(4, 5).Value = 0.000202
$ws.Cells.Item(5, 2).Value = 0.01594
$ws.Cells.Item(5, 3).Value = 47
$ws.Cells.Item(5, 4).Value = 216
$ws.Cells.Item(5, 5).Value = 0.000189
$ws.Cells.Item(6, 2).Value = 0.01464
$ws.Cells.Item(6, 3).Value = 40
$ws.Cells.Item(6, 4).Value = 178
$ws.Cells.Item(6, 5).Value = 0.000157
$ws.Cells.Item(7, 2).Value = 0.007308
$ws.Cells.Item(7, 3).Value = 39
$ws.Cells.Item(7, 4).Value = 171
$ws.Cells.Item(7, 5).Value = 0.000156
$ws.Cells.Item(8, 2).Value = 0.01355
$ws.Cells.Item(8, 3).Value = 37
$ws.Cells.Item(8, 4).Value = 175
$ws.Cells.Item(8, 5).Value = 0.00015
$ws.Cells.Item(9, 2).Value = 0.01034
$ws.Cells.Item(9, 3).Value = 144
$ws.Cells.Item(9, 4).Value = 755
$ws.Cells.Item(9, 5).Value = 0.000531
$ws.Cells.Item(10, 2).Value = 0.007681
$ws.Cells.Item(10, 3).Value = 707
$ws.Cells.Item(10, 4).Value = 3647
$ws.Cells.Item(10, 5).Value = 0.002213
$ws.Cells.Item(11, 2).Value = 0.005941
$ws.Cells.Item(11, 3).Value = 10
$ws.Cells.Item(11, 4).Value = 40
$ws.Cells.Item(11, 5).Value = 0.00006
